$d = $word.ActiveDocument

# --- Step 1: Replace paragraph 10 (the empty spacer paragraph right after the
# first details table) with the new "hasLinkedComplaints" begin paragraph,
# the new Linked Complaints table, and the new "hasLinkedComplaints" end
# paragraph, all in one shot by replacing the whole paragraph range.
$target = $d.Paragraphs(10)
$r = $target.Range

$fragment1 = @'
<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Segoe UI"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="en-CA"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Segoe UI"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="en-CA"/><w14:ligatures w14:val="none"/></w:rPr><w:t>{d.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Segoe UI"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="en-CA"/><w14:ligatures w14:val="none"/></w:rPr><w:t>ha</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Segoe UI"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="en-CA"/><w14:ligatures w14:val="none"/></w:rPr><w:t>sLinked</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Segoe UI"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="en-CA"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Complaints</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Segoe UI"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="en-CA"/><w14:ligatures w14:val="none"/></w:rPr><w:t>:ifEQ(true):showBegin}</w:t></w:r></w:p><w:tbl><w:tblPr><w:tblStyle w:val="PlainTable3"/><w:tblW w:w="10916" w:type="dxa"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="2552"/><w:gridCol w:w="8364"/></w:tblGrid><w:tr><w:trPr><w:cnfStyle w:val="100000000000" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:trHeight w:val="360"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000100" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="1" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="2552" w:type="dxa"/><w:noWrap/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:caps w:val="0"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Linked</w:t></w:r><w:commentRangeStart w:id="1"/><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:caps w:val="0"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w14:ligatures w14:val="none"/></w:rPr><w:t xml:space="preserve"> complaint(s)</w:t></w:r><w:commentRangeEnd w:id="1"/><w:r><w:rPr><w:rStyle w:val="CommentReference"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/></w:rPr><w:commentReference w:id="1"/></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="8364" w:type="dxa"/><w:noWrap/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:cnfStyle w:val="100000000000" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:tr><w:trPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:trHeight w:hRule="exact" w:val="20"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="10915" w:type="dxa"/><w:gridSpan w:val="2"/><w:tcBorders><w:right w:val="none" w:sz="0" w:space="0" w:color="auto"/></w:tcBorders><w:noWrap/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:spacing w:line="259" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:hRule="exact" w:val="20"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="10915" w:type="dxa"/><w:gridSpan w:val="2"/><w:tcBorders><w:right w:val="none" w:sz="0" w:space="0" w:color="auto"/></w:tcBorders><w:noWrap/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-CA"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:tr><w:trPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:trHeight w:val="360"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="10915" w:type="dxa"/><w:gridSpan w:val="2"/><w:tcBorders><w:right w:val="none" w:sz="0" w:space="0" w:color="auto"/></w:tcBorders><w:noWrap/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:spacing w:line="259" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>{d.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>linkedC</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>omplaints[i].id}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>{d</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>.linkedC</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>omplaints[i].</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>agency</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>{d</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>.linkedC</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>omplaints[i].</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>complaintTypeDescription</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>{d.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>linkedC</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>omplaints[i].</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>natureOfComplaint</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>}{d.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>linkedC</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>omplaints[i].</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>violationType</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>}{d.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>linkedC</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>omplaints[i].</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>girType</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>}</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="360"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="10915" w:type="dxa"/><w:gridSpan w:val="2"/><w:tcBorders><w:right w:val="none" w:sz="0" w:space="0" w:color="auto"/></w:tcBorders><w:noWrap/><w:vAlign w:val="center"/><w:hideMark/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>{d</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>.linkedC</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>omplaints[i</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>+1</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>].id}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>{d.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>linkedC</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>omplaints[i</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>+1</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>].</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>agency</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>{d</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>.linkedC</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>omplaints[i</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>+1</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>].</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>complaintTypeDescription</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>{d.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>linkedC</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>omplaints[i</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>+1</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>].</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>natureOfComplaint</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>}{d.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>linkedC</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>omplaints[i</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>+1</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>].</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>violationType</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>}{d.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>linkedC</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>omplaints[i</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>+1</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>].</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>girType</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>}</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Segoe UI"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="en-CA"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Segoe UI"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="en-CA"/><w14:ligatures w14:val="none"/></w:rPr><w:t>{d.hasLinked</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Segoe UI"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="en-CA"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Complaints</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="Times New Roman" w:hAnsi="BC Sans" w:cs="Segoe UI"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="en-CA"/><w14:ligatures w14:val="none"/></w:rPr><w:t>:showEnd}</w:t></w:r></w:p>
'@

$r.InsertXML($fragment1)

Write-Output "Step1 done"
